# Auto-generated: update FFXIV leve-profit market data cells (H:N) across 8 sheets
# per scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1042.64
$ws.Cells.Item(28, 9).Value = 554.4
$ws.Cells.Item(28, 10).Value = 1775
$ws.Cells.Item(28, 11).Value = 554.4
$ws.Cells.Item(28, 12).Value = 1775
$ws.Cells.Item(28, 13).Value = -69.39999999999998
$ws.Cells.Item(28, 14).Value = -2745
$ws.Cells.Item(74, 8).Value = 4418.4165
$ws.Cells.Item(74, 9).Value = 3752.1
$ws.Cells.Item(74, 11).Value = 3752.1
$ws.Cells.Item(74, 13).Value = -2816.1
$ws.Cells.Item(77, 8).Value = 4418.4165
$ws.Cells.Item(77, 9).Value = 3752.1
$ws.Cells.Item(77, 11).Value = 18760.5
$ws.Cells.Item(77, 13).Value = -14080.5
$ws.Cells.Item(112, 8).Value = 2945.6538
$ws.Cells.Item(112, 10).Value = 3304.3333
$ws.Cells.Item(112, 12).Value = 9912.999899999999
$ws.Cells.Item(112, 14).Value = -12128.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3646.6316
$ws.Cells.Item(2, 9).Value = 857.61536
$ws.Cells.Item(2, 10).Value = 9689.5
$ws.Cells.Item(2, 11).Value = 857.61536
$ws.Cells.Item(2, 12).Value = 9689.5
$ws.Cells.Item(2, 13).Value = -744.61536
$ws.Cells.Item(2, 14).Value = -9915.5
$ws.Cells.Item(32, 8).Value = 11123.374
$ws.Cells.Item(32, 9).Value = 4790.233
$ws.Cells.Item(32, 11).Value = 4790.233
$ws.Cells.Item(32, 13).Value = -4503.233
$ws.Cells.Item(74, 8).Value = 20530.217
$ws.Cells.Item(74, 9).Value = 3447.9285
$ws.Cells.Item(74, 11).Value = 3447.9285
$ws.Cells.Item(74, 13).Value = -2573.9285
$ws.Cells.Item(77, 8).Value = 20530.217
$ws.Cells.Item(77, 9).Value = 3447.9285
$ws.Cells.Item(77, 11).Value = 17239.6425
$ws.Cells.Item(77, 13).Value = -12871.6425
$ws.Cells.Item(86, 8).Value = 100314
$ws.Cells.Item(86, 10).Value = 100314
$ws.Cells.Item(86, 12).Value = 100314
$ws.Cells.Item(86, 14).Value = -102686
$ws.Cells.Item(89, 8).Value = 100314
$ws.Cells.Item(89, 10).Value = 100314
$ws.Cells.Item(89, 12).Value = 300942
$ws.Cells.Item(89, 14).Value = -312798
$ws.Cells.Item(116, 8).Value = 3646.6316
$ws.Cells.Item(116, 9).Value = 857.61536
$ws.Cells.Item(116, 10).Value = 9689.5
$ws.Cells.Item(116, 11).Value = 857.61536
$ws.Cells.Item(116, 12).Value = 9689.5
$ws.Cells.Item(116, 13).Value = 1436.38464
$ws.Cells.Item(116, 14).Value = -14277.5
$ws.Cells.Item(122, 8).Value = 3035.1191
$ws.Cells.Item(122, 9).Value = 2207.639
$ws.Cells.Item(122, 11).Value = 6622.917
$ws.Cells.Item(122, 13).Value = -4172.917
$ws.Cells.Item(132, 8).Value = 1394528.5
$ws.Cells.Item(132, 9).Value = 2236.1636
$ws.Cells.Item(132, 11).Value = 6708.4908
$ws.Cells.Item(132, 13).Value = -4178.4908
$ws.Cells.Item(139, 8).Value = 73755.336
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 73755.336
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 73755.336
$ws.Cells.Item(139, 13).ClearContents()
$ws.Cells.Item(139, 14).Value = -84035.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3646.6316
$ws.Cells.Item(3, 9).Value = 857.61536
$ws.Cells.Item(3, 10).Value = 9689.5
$ws.Cells.Item(3, 11).Value = 857.61536
$ws.Cells.Item(3, 12).Value = 9689.5
$ws.Cells.Item(3, 13).Value = -743.61536
$ws.Cells.Item(3, 14).Value = -9917.5
$ws.Cells.Item(20, 8).Value = 25190.73
$ws.Cells.Item(20, 9).Value = 10634.647
$ws.Cells.Item(20, 10).Value = 37563.4
$ws.Cells.Item(20, 11).Value = 10634.647
$ws.Cells.Item(20, 12).Value = 37563.4
$ws.Cells.Item(20, 13).Value = -10387.647
$ws.Cells.Item(20, 14).Value = -38057.4
$ws.Cells.Item(107, 8).Value = 3249.2163
$ws.Cells.Item(107, 9).Value = 3521.3547
$ws.Cells.Item(107, 11).Value = 3521.3547
$ws.Cells.Item(107, 13).Value = -1601.3547
$ws.Cells.Item(134, 8).Value = 10162.59
$ws.Cells.Item(134, 9).Value = 4855.724
$ws.Cells.Item(134, 11).Value = 14567.172
$ws.Cells.Item(134, 13).Value = -12032.172
$ws.Cells.Item(137, 8).Value = 89500
$ws.Cells.Item(137, 10).Value = 89500
$ws.Cells.Item(137, 12).Value = 89500
$ws.Cells.Item(137, 14).Value = -99700

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 2266.1714
$ws.Cells.Item(107, 9).Value = 780.11536
$ws.Cells.Item(107, 10).Value = 6559.222
$ws.Cells.Item(107, 11).Value = 780.11536
$ws.Cells.Item(107, 12).Value = 6559.222
$ws.Cells.Item(107, 13).Value = 1139.88464
$ws.Cells.Item(107, 14).Value = -10399.222
$ws.Cells.Item(124, 8).Value = 28326
$ws.Cells.Item(124, 10).Value = 28326
$ws.Cells.Item(124, 12).Value = 28326
$ws.Cells.Item(124, 14).Value = -33236
$ws.Cells.Item(134, 8).Value = 47629212
$ws.Cells.Item(134, 9).Value = 2703
$ws.Cells.Item(134, 10).Value = 71442460
$ws.Cells.Item(134, 11).Value = 8109
$ws.Cells.Item(134, 12).Value = 214327380
$ws.Cells.Item(134, 13).Value = -5574
$ws.Cells.Item(134, 14).Value = -214332450
$ws.Cells.Item(138, 8).Value = 77308.78
$ws.Cells.Item(138, 10).Value = 77308.78
$ws.Cells.Item(138, 12).Value = 77308.78
$ws.Cells.Item(138, 14).Value = -87588.78
$ws.Cells.Item(141, 8).Value = 149447.6
$ws.Cells.Item(141, 10).Value = 156269.62
$ws.Cells.Item(141, 12).Value = 156269.62
$ws.Cells.Item(141, 14).Value = -166629.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 559
$ws.Cells.Item(11, 9).Value = 559
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 1677
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -1537
$ws.Cells.Item(11, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 1839038
$ws.Cells.Item(107, 10).Value = 2404766
$ws.Cells.Item(107, 12).Value = 7214298
$ws.Cells.Item(107, 14).Value = -7218138
$ws.Cells.Item(117, 8).Value = 2827
$ws.Cells.Item(117, 10).Value = 3608.3845
$ws.Cells.Item(117, 12).Value = 10825.1535
$ws.Cells.Item(117, 14).Value = -17709.1535
$ws.Cells.Item(131, 8).Value = 1432.1786
$ws.Cells.Item(131, 9).Value = 738.5714
$ws.Cells.Item(131, 10).Value = 1495.2338
$ws.Cells.Item(131, 11).Value = 2215.7142
$ws.Cells.Item(131, 12).Value = 4485.7014
$ws.Cells.Item(131, 13).Value = 2824.2858
$ws.Cells.Item(131, 14).Value = -14565.7014
$ws.Cells.Item(136, 8).Value = 2520.5715
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 14501
$ws.Cells.Item(70, 9).Value = 4500
$ws.Cells.Item(70, 10).Value = 17001.25
$ws.Cells.Item(70, 11).Value = 4500
$ws.Cells.Item(70, 12).Value = 17001.25
$ws.Cells.Item(70, 13).Value = -4230
$ws.Cells.Item(70, 14).Value = -17541.25
$ws.Cells.Item(73, 8).Value = 14501
$ws.Cells.Item(73, 9).Value = 4500
$ws.Cells.Item(73, 10).Value = 17001.25
$ws.Cells.Item(73, 11).Value = 4500
$ws.Cells.Item(73, 12).Value = 17001.25
$ws.Cells.Item(73, 13).Value = -3564
$ws.Cells.Item(73, 14).Value = -18873.25
$ws.Cells.Item(80, 8).Value = 8617.444
$ws.Cells.Item(80, 9).Value = 6988.8887
$ws.Cells.Item(80, 10).Value = 9703.147999999999
$ws.Cells.Item(80, 11).Value = 6988.8887
$ws.Cells.Item(80, 12).Value = 9703.147999999999
$ws.Cells.Item(80, 13).Value = -5990.8887
$ws.Cells.Item(80, 14).Value = -11699.148
$ws.Cells.Item(83, 8).Value = 8617.444
$ws.Cells.Item(83, 9).Value = 6988.8887
$ws.Cells.Item(83, 10).Value = 9703.147999999999
$ws.Cells.Item(83, 11).Value = 34944.4435
$ws.Cells.Item(83, 12).Value = 48515.74
$ws.Cells.Item(83, 13).Value = -29952.4435
$ws.Cells.Item(83, 14).Value = -58499.74
$ws.Cells.Item(107, 8).Value = 1253.1818
$ws.Cells.Item(107, 9).Value = 797.6
$ws.Cells.Item(107, 11).Value = 797.6
$ws.Cells.Item(107, 13).Value = 1122.4
$ws.Cells.Item(113, 8).Value = 1233.1052
$ws.Cells.Item(113, 9).Value = 1040.7693
$ws.Cells.Item(113, 10).Value = 1649.8334
$ws.Cells.Item(113, 11).Value = 1040.7693
$ws.Cells.Item(113, 12).Value = 1649.8334
$ws.Cells.Item(113, 13).Value = 1129.2307
$ws.Cells.Item(113, 14).Value = -5989.8334
$ws.Cells.Item(126, 8).Value = 6602.5186
$ws.Cells.Item(126, 9).Value = 3921.111
$ws.Cells.Item(126, 11).Value = 11763.333
$ws.Cells.Item(126, 13).Value = -9293.332999999999
$ws.Cells.Item(132, 8).Value = 4570.654
$ws.Cells.Item(132, 9).Value = 1652.4595
$ws.Cells.Item(132, 10).Value = 11768.866
$ws.Cells.Item(132, 11).Value = 4957.3785
$ws.Cells.Item(132, 12).Value = 35306.598
$ws.Cells.Item(132, 13).Value = -2427.3785
$ws.Cells.Item(132, 14).Value = -40366.598

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 21999.5
$ws.Cells.Item(22, 9).Value = 2999
$ws.Cells.Item(22, 10).Value = 28333
$ws.Cells.Item(22, 11).Value = 2999
$ws.Cells.Item(22, 12).Value = 28333
$ws.Cells.Item(22, 13).Value = -2704
$ws.Cells.Item(22, 14).Value = -28923
$ws.Cells.Item(27, 8).Value = 21999.5
$ws.Cells.Item(27, 9).Value = 2999
$ws.Cells.Item(27, 10).Value = 28333
$ws.Cells.Item(27, 11).Value = 2999
$ws.Cells.Item(27, 12).Value = 28333
$ws.Cells.Item(27, 13).Value = -2892
$ws.Cells.Item(27, 14).Value = -28547
$ws.Cells.Item(61, 8).Value = 2750.8215
$ws.Cells.Item(61, 9).Value = 1717.762
$ws.Cells.Item(61, 10).Value = 5850
$ws.Cells.Item(61, 11).Value = 1717.762
$ws.Cells.Item(61, 12).Value = 5850
$ws.Cells.Item(61, 13).Value = -1515.762
$ws.Cells.Item(61, 14).Value = -6254
$ws.Cells.Item(63, 8).Value = 20333.334
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 20333.334
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 20333.334
$ws.Cells.Item(63, 13).ClearContents()
$ws.Cells.Item(63, 14).Value = -21831.334
$ws.Cells.Item(66, 8).Value = 20333.334
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 20333.334
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 61000.00199999999
$ws.Cells.Item(66, 13).ClearContents()
$ws.Cells.Item(66, 14).Value = -68488.00199999999
$ws.Cells.Item(113, 8).Value = 2750.8215
$ws.Cells.Item(113, 9).Value = 1717.762
$ws.Cells.Item(113, 10).Value = 5850
$ws.Cells.Item(113, 11).Value = 1717.762
$ws.Cells.Item(113, 12).Value = 5850
$ws.Cells.Item(113, 13).Value = 452.2380000000001
$ws.Cells.Item(113, 14).Value = -10190
$ws.Cells.Item(133, 8).Value = 79991
$ws.Cells.Item(133, 10).Value = 79991
$ws.Cells.Item(133, 12).Value = 79991
$ws.Cells.Item(133, 14).Value = -85051

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 4005
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(96, 8).Value = 2363.4546
$ws.Cells.Item(96, 9).Value = 1150
$ws.Cells.Item(96, 10).Value = 2633.111
$ws.Cells.Item(96, 11).Value = 1150
$ws.Cells.Item(96, 12).Value = 2633.111
$ws.Cells.Item(96, 13).Value = 223
$ws.Cells.Item(96, 14).Value = -5379.111
$ws.Cells.Item(98, 8).Value = 59997.668
$ws.Cells.Item(98, 10).Value = 59998
$ws.Cells.Item(98, 12).Value = 59998
$ws.Cells.Item(98, 14).Value = -65988
$ws.Cells.Item(136, 8).Value = 20170.646
$ws.Cells.Item(136, 9).Value = 3571
$ws.Cells.Item(136, 11).Value = 10713
$ws.Cells.Item(136, 13).Value = -8163
